# Updates loading_percent values for rows 2-25 (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","G","H","J","M","N","O")

$data = @(
    @(13.24945670553124,6.931691617177907,14.53418642754317,15.73022987164546,31.77039380191395,14.94847978023523,9.221411298107018,17.3298259617475,17.70740358767744,23.17966220647659),
    @(12.78469178007335,6.514330546989612,14.51853349932995,15.74717638812273,31.73654461844281,14.98965538811909,9.242858625980887,17.19200000480822,17.76131818696163,23.22690651297395),
    @(12.49229289391237,6.242824504779835,14.51192525340317,15.76049987877704,31.7267528106048,15.01770552241795,9.256907866277395,17.10964036018673,17.79623308816923,23.26159497780436),
    @(12.37155408816552,6.128373176147104,14.50998993690901,15.7666628728418,31.7255267115056,15.02983133730265,9.262854836129467,17.07667767393577,17.81091765257753,23.27715544443044),
    @(12.35141555796186,6.109139298003861,14.50971440007858,15.76773052726407,31.72549001377168,15.03188677772154,9.263855732844503,17.07124130399428,17.81338361307305,23.27982517613977),
    @(12.49067072526217,6.241296366693295,14.51189608259482,15.7605800253892,31.72672508508296,15.017866241993,9.256987170662709,17.10919334787274,17.79642927977709,23.26179906863571),
    @(13.09077033324556,6.790961166579973,14.52816776056279,15.735467247391,31.75644226214754,14.96210206564654,9.228623820352407,17.28185138214938,17.72561806706177,23.19477098651298),
    @(14.20454493484585,7.74721714579625,14.58376867579047,15.70938109779845,31.90178996949284,14.8747521748539,9.179972486580514,17.63706424893988,17.60108075567433,23.10854962499131),
    @(14.9759075274948,8.375091116827914,14.6388526431153,15.70433121979547,32.06127207697141,14.8240419817075,9.148453178649016,17.90631980182479,17.51824934206643,23.07295155968842),
    @(15.31514709252121,8.644475139881342,14.66694700039253,15.70509382414728,32.14512075129784,14.80390659432912,9.135026639148213,18.03020284075776,17.48243526676477,23.06281185375862),
    @(15.44182769476226,8.74414883308606,14.67801665040653,15.70582172271298,32.17847978429269,14.79670432573237,9.130073096779041,18.07728046837286,17.46914079873826,23.05984432373402),
    @(15.41462556350353,8.722786216388858,14.6756135338914,15.70564544352069,32.17122414179639,14.79823665721873,9.131134118453964,18.06713464414305,17.47199211419721,23.06044462519634),
    @(15.32560545638383,8.652722173666486,14.66784910523596,15.70514491415857,32.14783311582118,14.80330558610528,9.134616488644646,18.03407282377475,17.48133616557514,23.06255022668682),
    @(15.27084307853034,8.609501689604116,14.66314910985219,15.70489548113558,32.13371419310506,14.80646550486215,9.136766565358721,18.01384205419836,17.48709448474142,23.06395358733205),
    @(14.95349350770363,8.357159327422261,14.63707723502574,15.70434288339625,32.05601831588091,14.82541696587131,9.149348954569282,17.89824887075828,17.5206273600583,23.07373618405859),
    @(14.75574464270238,8.198194802401696,14.6218567831963,15.70478702562093,32.01123697938748,14.83779486883489,9.157301140692931,17.82766912398725,17.54167610958964,23.08128937197279),
    @(14.64091294018505,8.10523402390748,14.6133885552446,15.70533055722778,31.98654473654445,14.84519032848961,9.161960864546645,17.78720715611181,17.55395851171429,23.08620353721628),
    @(14.60184882482836,8.073496633699468,14.61057067455327,15.70556408734331,31.9783677315074,14.8477416843854,9.163553319112287,17.77353145009333,17.55814732757712,23.08796519087189),
    @(14.77690915077335,8.215275065043553,14.62344744807141,15.70470993708845,32.01589394378112,14.83644864705499,9.156445735399547,17.83516888196758,17.53941725420432,23.08042633980791),
    @(15.35180196239125,8.67336508255168,14.67011805890918,15.70528002249162,32.15466015790989,14.80180524437644,9.13359008478163,18.04377966194361,17.47858433662106,23.06190808029304),
    @(15.7171003843239,8.959137236163164,14.703128766775,15.70821153413186,32.2547106447409,14.7816271407782,9.119414838134237,18.18106904522764,17.44038559446144,23.05488929144638),
    @(15.52311945120581,8.807860945543542,14.68528279425234,15.70641314331269,32.20046202838359,14.79217092623668,9.12691079063781,18.10771978227553,17.46063059527936,23.05816977693841),
    @(14.76734422944523,8.207557968312425,14.6227274289986,15.70474389107266,32.01378524838543,14.8370564041518,9.156832190284742,17.83177787964695,17.5404379182307,23.08081473582297),
    @(13.91092716345117,7.501618009458705,14.56621023425148,15.71395811046676,31.85317565830668,14.89602142528526,9.192390342024442,17.5393786873996,17.63324486802708,23.12701310233107)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $cellRef = "$($cols[$c])$rowNum"
        $ws.Range($cellRef).Value = $row[$c]
    }
}
